$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the existing table (Table1) that lists the by-element export properties
$tbl = $ws.ListObjects.Item(1)

# Add a new row to the table for the "storagebus" element with the "puVmagAngle" property
$newRow = $tbl.ListRows.Add()
$newRow.Range.Item(1, 1).Value = "storagebus"
$newRow.Range.Item(1, 2).Value = "puVmagAngle"

# Match formatting of the new row to the last existing data row's 3rd (Property 3) column style
$ws.Range("C4").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)

# Widen column A (Element) to fit the updated content / per author's manual resize
$ws.Columns("A").ColumnWidth = 32.42

# Move the active selection, matching the saved cursor position in the workbook
$ws.Range("B9").Select()
